$wb = $excel.ActiveWorkbook

$sheet2 = $wb.Worksheets.Item("moving_average")
$sheet2.Range("D18").Value = 156
$sheet2.Range("D19").Value = 156
$sheet2.Range("D20").Value = 156
$sheet2.Range("D21").Value = 156

$sheet3 = $wb.Worksheets.Item("simple_exponential_smoothing")
$sheet3.Range("D19").Value = 139.4751512085888
$sheet3.Range("D20").Value = 139.4751512085888
$sheet3.Range("D21").Value = 139.4751512085888
$sheet3.Range("D22").Value = 139.4751512085888

$sheet5 = $wb.Worksheets.Item("winter_trendseason")
$sheet5.Range("C3").Value = 108.7856364997036
$sheet5.Range("D3").Value = 2.875737309644022
$sheet5.Range("E3").Value = 0.9861455576906619
$sheet5.Range("C4").Value = 111.7132263929941
$sheet5.Range("D4").Value = 2.88092256800867
$sheet5.Range("E4").Value = 0.9448690013829411
$sheet5.Range("F4").Value = 105.4558787096639
$sheet5.Range("G4").Value = -0.5441212903360508
$sheet5.Range("H4").Value = 0.5441212903360508
$sheet5.Range("I4").Value = 64.32144739385711
$sheet5.Range("J4").Value = 5.936573572355968
$sheet5.Range("K4").Value = 0.5133219720151423
$sheet5.Range("L4").Value = 6.036776217832002
$sheet5.Range("M4").Value = 1.816688436956366
$sheet5.Range("C5").Value = 113.7925059591224
$sheet5.Range("D5").Value = 2.800758267820626
$sheet5.Range("E5").Value = 1.024000894415179
$sheet5.Range("F5").Value = 118.1863588276995
$sheet5.Range("G5").Value = 9.186358827699536
$sheet5.Range("H5").Value = 9.186358827699536
$sheet5.Range("I5").Value = 71.01069443298914
$sheet5.Range("J5").Value = 7.019835324137158
$sheet5.Range("K5").Value = 8.427852135504162
$sheet5.Range("L5").Value = 6.833801523722722
$sheet5.Range("M5").Value = 2.844976052795674
$sheet5.Range("C6").Value = 117.7304273612062
$sheet5.Range("D6").Value = 2.914474581246947
$sheet5.Range("E6").Value = 1.039236854531764
$sheet5.Range("F6").Value = 119.9960915916242
$sheet5.Range("G6").Value = -13.00390840837579
$sheet5.Range("H6").Value = 13.00390840837579
$sheet5.Range("I6").Value = 95.53342929809853
$sheet5.Range("J6").Value = 8.515853595196816
$sheet5.Range("K6").Value = 9.777374743139694
$sheet5.Range("L6").Value = 7.569694828576965
$sheet5.Range("M6").Value = 0.8181628424534405
$sheet5.Range("C7").Value = 121.6512350567574
$sheet5.Range("D7").Value = 3.015107892677377
$sheet5.Range("E7").Value = 0.9943938709581061
$sheet5.Range("F7").Value = 118.9734341085757
$sheet5.Range("G7").Value = -11.02656589142433
$sheet5.Range("H7").Value = 11.02656589142433
$sheet5.Range("I7").Value = 100.7437745100633
$sheet5.Range("J7").Value = 9.017996054442317
$sheet5.Range("K7").Value = 8.481973762634096
$sheet5.Range("L7").Value = 7.752150615388392
$sheet5.Range("M7").Value = -0.4501233847913648
$sheet5.Range("C8").Value = 124.495522796645
$sheet5.Range("D8").Value = 2.998025877398394
$sheet5.Range("E8").Value = 0.9435581427553752
$sheet5.Range("F8").Value = 117.7933629686957
$sheet5.Range("G8").Value = 1.793362968695732
$sheet5.Range("H8").Value = 1.793362968695732
$sheet5.Range("I8").Value = 84.4891705479676
$sheet5.Range("J8").Value = 7.81389054015122
$sheet5.Range("K8").Value = 1.546002559220458
$sheet5.Range("L8").Value = 6.717792606027071
$sheet5.Range("M8").Value = -0.289976923495676
$sheet5.Range("C9").Value = 127.7085722082084
$sheet5.Range("D9").Value = 3.019528230814897
$sheet5.Range("E9").Value = 1.025744166456996
$sheet5.Range("F9").Value = 130.5535078743856
$sheet5.Range("G9").Value = -2.446492125614384
$sheet5.Range("H9").Value = 2.446492125614384
$sheet5.Range("I9").Value = 73.27433528692839
$sheet5.Range("J9").Value = 7.04711933807453
$sheet5.Range("K9").Value = 1.839467763619838
$sheet5.Range("L9").Value = 6.020889057111751
$sheet5.Range("M9").Value = -0.6686902603620367
$sheet5.Range("C10").Value = 129.0084045072139
$sheet5.Range("D10").Value = 2.847558637633959
$sheet5.Range("E10").Value = 1.025229791521208
$sheet5.Range("F10").Value = 135.8574598991631
$sheet5.Range("G10").Value = 19.85745989916308
$sheet5.Range("H10").Value = 19.85745989916308
$sheet5.Range("I10").Value = 113.4048825819211
$sheet5.Range("J10").Value = 8.648411908210598
$sheet5.Range("K10").Value = 17.11849991307162
$sheet5.Range("L10").Value = 7.408090414106733
$sheet5.Range("M10").Value = 1.751202416689387
$sheet5.Range("C11").Value = 132.4789471294101
$sheet5.Range("D11").Value = 2.909857036090178
$sheet5.Range("E11").Value = 0.9991219783889517
$sheet5.Range("F11").Value = 131.1167616005146
$sheet5.Range("G11").Value = -6.88323839948535
$sheet5.Range("H11").Value = 6.88323839948535
$sheet5.Range("I11").Value = 106.0686701688354
$sheet5.Range("J11").Value = 8.452281518352237
$sheet5.Range("K11").Value = 4.987853912670544
$sheet5.Range("L11").Value = 7.139175247280489
$sheet5.Range("M11").Value = 0.9774735279178171
$sheet5.Range("C12").Value = 135.6036835842328
$sheet5.Range("D12").Value = 2.931344977963434
$sheet5.Range("E12").Value = 0.9450699307188639
$sheet5.Range("F12").Value = 127.7472086082706
$sheet5.Range("G12").Value = -2.252791391729374
$sheet5.Range("H12").Value = 2.252791391729374
$sheet5.Range("I12").Value = 95.96931005741681
$sheet5.Range("J12").Value = 7.83233250568995
$sheet5.Range("K12").Value = 1.732916455176442
$sheet5.Range("L12").Value = 6.598549368070084
$sheet5.Range("M12").Value = 0.7672159013427391
$sheet5.Range("C13").Value = 138.9648289438416
$sheet5.Range("D13").Value = 2.974325016127965
$sheet5.Range("E13").Value = 1.028951911468478
$sheet5.Range("F13").Value = 142.1014973976262
$sheet5.Range("G13").Value = -4.898502602373838
$sheet5.Range("H13").Value = 4.898502602373838
$sheet5.Range("I13").Value = 89.42622075633011
$sheet5.Range("J13").Value = 7.565620696297577
$sheet5.Range("K13").Value = 3.332314695492407
$sheet5.Range("L13").Value = 6.301618943290296
$sheet5.Range("M13").Value = 0.1467939624753597
$sheet5.Range("C14").Value = 141.5423429879891
$sheet5.Range("D14").Value = 2.934643918929919
$sheet5.Range("E14").Value = 1.022323645765695
$sheet5.Range("F14").Value = 145.5202492230762
$sheet5.Range("G14").Value = 4.520249223076178
$sheet5.Range("H14").Value = 4.520249223076178
$sheet5.Range("I14").Value = 83.67675677986267
$sheet5.Range("J14").Value = 7.311839740195794
$sheet5.Range("K14").Value = 3.205850512819984
$sheet5.Range("L14").Value = 6.043638240751103
$sheet5.Range("M14").Value = 0.7700984791442528
$sheet5.Range("C15").Value = 144.4454472453041
$sheet5.Range("D15").Value = 2.931489952768429
$sheet5.Range("E15").Value = 0.998901396136564
$sheet5.Range("F15").Value = 144.3501329901156
$sheet5.Range("G15").Value = 0.3501329901155543
$sheet5.Range("H15").Value = 0.3501329901155543
$sheet5.Range("I15").Value = 77.24951342070149
$sheet5.Range("J15").Value = 6.776323836343467
$sheet5.Range("K15").Value = 0.2431479098024683
$sheet5.Range("L15").Value = 5.597446676831978
$sheet5.Range("M15").Value = 0.8826274833131647
$sheet5.Range("C16").Value = 147.6358216759156
$sheet5.Range("D16").Value = 2.957378400552735
$sheet5.Range("E16").Value = 0.9467455567349878
$sheet5.Range("F16").Value = 139.2815118273408
$sheet5.Range("G16").Value = -2.718488172659249
$sheet5.Range("H16").Value = 2.718488172659249
$sheet5.Range("I16").Value = 72.25956088671482
$sheet5.Range("J16").Value = 6.486478431794595
$sheet5.Range("K16").Value = 1.914428290605105
$sheet5.Range("L16").Value = 5.334373934958629
$sheet5.Range("M16").Value = 0.5029665195733909
$sheet5.Range("C17").Value = 151.4719734120989
$sheet5.Range("D17").Value = 3.0452557341158
$sheet5.Range("E17").Value = 1.03498776299782
$sheet5.Range("F17").Value = 154.953161072837
$sheet5.Range("G17").Value = -10.04683892716301
$sheet5.Range("H17").Value = 10.04683892716301
$sheet5.Range("I17").Value = 74.17152165615771
$sheet5.Range("J17").Value = 6.723835798152489
$sheet5.Range("K17").Value = 6.088993289189702
$sheet5.Range("L17").Value = 5.384681891907367
$sheet5.Range("M17").Value = -1.009001059767039
$sheet5.Range("C18").Value = 155.8406891605757
$sheet5.Range("D18").Value = 3.177601735551892
$sheet5.Range("E18").Value = 1.03110208390814
$sheet5.Range("F18").Value = 157.9666170343716
$sheet5.Range("G18").Value = -15.03338296562845
$sheet5.Range("H18").Value = 15.03338296562845
$sheet5.Range("I18").Value = 83.66096426460081
$sheet5.Range("J18").Value = 7.243182496119736
$sheet5.Range("K18").Value = 8.689816743137831
$sheet5.Range("L18").Value = 5.591252820109272
$sheet5.Range("M18").Value = -3.012175990781939
$sheet5.Range("E19").Value = 0.998901396136564
$sheet5.Range("F19").Value = 151.3055657169062
$sheet5.Range("E20").Value = 0.9467455567349878
$sheet5.Range("F20").Value = 146.2885001331808
$sheet5.Range("E21").Value = 1.03498776299782
$sheet5.Range("F21").Value = 163.0752437586711
$sheet5.Range("E22").Value = 1.03110208390814
$sheet5.Range("F22").Value = 165.6029760393337
